$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '51.038.34'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -1.88%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.894.29'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -1.53%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.996'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.78%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '367.93'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +3.80%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '101.86'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -5.71%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.537'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -4.98%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.997'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.51%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.587'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -5.55%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '36.64'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -4.32%  '
$ws.Range('E11').Value = '  +1.25%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0832'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -3.82%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '18.31'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -4.39%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '3.342.31'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -2.74%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.37'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -4.98%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '2.882.82'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -2.30%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.929'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -4.31%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '50.703.92'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -2.63%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '3.27'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -5.45%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.19'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -4.40%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '12.84'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -5.77%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.0₃0941'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -3.41%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '67.87'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.97%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '258.04'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -2.84%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.67'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -3.39%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.171'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -2.54%  '
$ws.Range('E27').Value = '  +0.16%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '25.55'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -4.89%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '6.95'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -8.18%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.103'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -2.69%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '9.84'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -4.54%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '5.99'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.59%  '
$ws.Range('E33').Value = '  -2.28%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '34.45'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -6.17%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '50.71'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -2.77%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.997'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.40%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.0414'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -5.23%  '
$ws.Range('B38').Value = 'EnergySwap'
$ws.Range('C38').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '25.84'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +12.70%  '
$ws.Range('B39').Value = 'LidoDAOToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.05'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -4.04%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.63'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -2.34%  '
$ws.Range('B41').Value = 'Celestia'
$ws.Range('C41').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '16.90'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -5.48%  '
$ws.Range('B42').Value = 'ARBITRUM'
$ws.Range('C42').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.84'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -7.93%  '
$ws.Range('B43').Value = 'Stellar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.112'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -4.53%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '118.00'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.28%  '
$ws.Range('E45').Value = '  -3.57%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.022.32'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -4.84%  '
$ws.Range('E47').Value = '  -6.17%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '3.15'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -6.99%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '3.250.63'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -1.54%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.234'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -2.44%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0310'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -10.04%  '
